# Regenerate Report for Handback.
#
# The "Correspond Handoff Datetime" (col E) and "Correspond Handback
# DateTime" (col H) values for the 525d3822-... record are refreshed on
# both the zh-cn and de-de sheets. In the source workbook the row-4 and
# row-5 records for these two columns happen to share the same text
# value, so both rows are (re)written with the new timestamp to keep
# the worksheet content consistent with that shared value.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-13 18:23:05"
$wsZhCn.Range("E5").Value = "2016-03-13 18:23:05"
$wsZhCn.Range("H4").Value = "2016-03-13 18:23:25"
$wsZhCn.Range("H5").Value = "2016-03-13 18:23:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-13 18:23:09"
$wsDeDe.Range("E5").Value = "2016-03-13 18:23:09"
$wsDeDe.Range("H4").Value = "2016-03-13 18:23:31"
$wsDeDe.Range("H5").Value = "2016-03-13 18:23:31"
